$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 3 data
$ws.Range("B3").Value = "Amiraaa"
$ws.Range("C3").Value = "1116332215&"

# Match the number format used by C2 (custom format 0;[Red]0)
$ws.Range("C3").NumberFormat = $ws.Range("C2").NumberFormat

# Update the active selection to C3, as in the edited workbook
$ws.Range("C3").Select()
